# Predicted_LaLiga2025_26_table_matchday_4.xlsx
# Add extra prediction columns (WIN / TOP4 / TOP5 / TOP6 / RELEGATION) between
# Team and ExpPoints, shifting ExpPoints to column H, refreshing the ExpPoints
# values, and re-ordering the few teams whose predicted rank moved - in
# preparation for a Monte Carlo simulation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------------
# Old layout: A=Rank B=Team C=ExpPoints
# New layout: A=Rank B=Team C=WIN D=TOP4 E=TOP5 F=TOP6 G=RELEGATION H=ExpPoints
$ws.Cells.Item(1, 8).Value = "ExpPoints"   # H1 (move the header out of C1 first)
$ws.Cells.Item(1, 3).Value = "WIN"         # C1
$ws.Cells.Item(1, 4).Value = "TOP4"        # D1
$ws.Cells.Item(1, 5).Value = "TOP5"        # E1
$ws.Cells.Item(1, 6).Value = "TOP6"        # F1
$ws.Cells.Item(1, 7).Value = "RELEGATION"  # G1

# copy the bold/centered/bordered header style from A1 onto the new header cells
$ws.Range("A1").Copy()
$ws.Range("C1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows ----------------------------------------------------------------
# Row -> Team, ExpPoints (refreshed values / refreshed team order)
$rows = @(
    @(2,  "Barcelona",          85.70803318802812),
    @(3,  "Real Madrid",        84.79025131044398),
    @(4,  "Atlético de Madrid", 71.13042728157198),
    @(5,  "Villarreal",         64.34317872761542),
    @(6,  "Real Betis",         61.63174444565892),
    @(7,  "Athletic Club",      56.03483543184727),
    @(8,  "Celta de Vigo",      49.72517546778344),
    @(9,  "Espanyol",           49.67067174980772),
    @(10, "Rayo Vallecano",     48.97723237696417),
    @(11, "Getafe",             47.66788018967382),
    @(12, "Osasuna",            46.22750982726023),
    @(13, "Sevilla",            45.43281201888133),
    @(14, "Real Sociedad",      44.47843563682105),
    @(15, "Valencia",           44.34502795488061),
    @(16, "Alavés",             42.6308513766258),
    @(17, "Elche",              42.42812344465041),
    @(18, "Mallorca",           40.94283105730698),
    @(19, "Levante",            35.54252111563662),
    @(20, "Girona",             32.70439630721698),
    @(21, "Real Oviedo",        32.6010896571655)
)

foreach ($r in $rows) {
    $row = $r[0]
    $team = $r[1]
    $expPoints = $r[2]

    $ws.Cells.Item($row, 8).Value = $expPoints      # H: ExpPoints (write before clearing C)
    $ws.Cells.Item($row, 2).Value = $team           # B: Team

    # leading apostrophe forces a (blank) text entry instead of Excel deleting
    # the cell outright, keeping a present-but-empty placeholder cell, like
    # the source workbook's empty inlineStr cells for the new WIN/TOP4/TOP5/
    # TOP6/RELEGATION columns (still to be filled in by the Monte Carlo sim).
    $placeholders = $ws.Range($ws.Cells.Item($row, 3), $ws.Cells.Item($row, 7))
    $placeholders.Value = "'"
    $placeholders.Style = "Normal"   # drop the auto-added quote-prefix style
}

$ws.Range("A1").Select() | Out-Null
